# Apply the crypto price/volume refresh for this run (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for existing coin rows.
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of silently converting them to numbers (and dropping trailing zeros).
$ws.Range("D2").Value = "24.406.26"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.652.28"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'312.19"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").Value = "'0.3917"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "'50.68"
$ws.Range("E10").Value = "  -5.42%  "
$ws.Range("D11").Value = "'1.377"
$ws.Range("E11").Value = "  -6.35%  "
$ws.Range("D12").Value = "'0.08568"
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("D13").Value = "'25.05"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").Value = "'7.203"
$ws.Range("E14").Value = "  -4.10%  "
$ws.Range("D15").Value = "'0.00001305"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "'7.607"
$ws.Range("E16").Value = "  -4.76%  "
$ws.Range("D17").Value = "1.654.25"
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("D18").Value = "'93.28"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("D19").Value = "'0.06954"
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("D20").Value = "'21.08"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'7.005"
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("D24").Value = "24.406.99"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'2.348"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'2.777"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "'22.68"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").Value = "'158.74"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").Value = "'5.734"
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("D30").Value = "'145.12"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "'8.125"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "'2.510"
$ws.Range("E32").Value = "  +10.61%  "
$ws.Range("D33").Value = "1.845.67"
$ws.Range("E33").Value = "  -3.04%  "

# Source ranking reshuffled rows 34-36: VeChain moved up to 34, ImmutableX
# dropped to 35, and Hedera dropped to 36 (each keeps its own updated data).
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.03024"
$ws.Range("E34").Value = "  -5.47%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.009"
$ws.Range("E35").Value = "  -2.21%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.08149"
$ws.Range("E36").Value = "  -5.45%  "

$ws.Range("D37").Value = "'6.848"
$ws.Range("E37").Value = "  -5.84%  "
$ws.Range("D38").Value = "'0.2761"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").Value = "'0.09518"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").Value = "'1.496"
$ws.Range("E40").Value = "  +1.34%  "
$ws.Range("D41").Value = "'10.21"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").Value = "'0.7774"
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("D43").Value = "'13.32"
$ws.Range("E43").Value = "  -6.09%  "
$ws.Range("D44").Value = "'16.16"
$ws.Range("E44").Value = "  -7.68%  "
$ws.Range("D45").Value = "'2.558"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("D46").Value = "'0.7006"
$ws.Range("E46").Value = "  -5.73%  "
$ws.Range("D47").Value = "'4.149"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "'0.08546"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'1.300"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("D51").Value = "'136.51"
$ws.Range("E51").Value = "  -2.74%  "
